$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status values: Github/Environment -> OK (1), VSCode -> OK (1)
$ws.Range("B2").Value = 1
$ws.Range("B4").Value = 1

# Update the current selection / view (topLeftCell reset, selection moved to D6)
$ws.Range("D6").Select()
